$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D, E changed
$ws.Range("D2").Value = "34.063.06"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3: D, E changed
$ws.Range("D3").Value = "1.779.93"
$ws.Range("E3").Value = "  -0.51%  "

# Row 4: E changed
$ws.Range("E4").Value = "  +0.35%  "

# Row 5: D, E changed
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.19"
$ws.Range("E5").Value = "  -0.76%  "

# Row 6: E changed
$ws.Range("E6").Value = "  +0.00%  "

# Row 7: E changed
$ws.Range("E7").Value = "  +0.29%  "

# Row 8: D, E changed
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.64"
$ws.Range("E8").Value = "  -1.99%  "

# Row 9: E changed
$ws.Range("E9").Value = "  -1.75%  "

# Row 10: D, E changed
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0684"
$ws.Range("E10").Value = "  -0.09%  "

# Row 11: D, E changed
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("E11").Value = "  +0.79%  "

# Row 12: D, E changed
$ws.Range("D12").Value = "2.036.54"
$ws.Range("E12").Value = "  -0.52%  "

# Row 13: D, E changed
$ws.Range("D13").Value = "1.780.16"
$ws.Range("E13").Value = "  -0.27%  "

# Row 14: D, E changed
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.87"
$ws.Range("E14").Value = "  -4.15%  "

# Row 15: D, E changed
$ws.Range("D15").Value = "34.074.27"
$ws.Range("E15").Value = "  +0.02%  "

# Row 16: D, E changed
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.619"
$ws.Range("E16").Value = "  -0.88%  "

# Row 17: E changed
$ws.Range("E17").Value = "  -0.50%  "

# Row 18: D, E changed
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.49"
$ws.Range("E18").Value = "  -0.72%  "

# Row 19: D, E changed
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.15"
$ws.Range("E19").Value = "  +1.12%  "

# Row 20: D, E changed
$ws.Range("D20").Value = "0.0₃0785"
$ws.Range("E20").Value = "  +1.38%  "

# Row 21: E changed
$ws.Range("E21").Value = "  +0.42%  "

# Row 22: D, E changed
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.84"
$ws.Range("E22").Value = "  +1.02%  "

# Row 23: D, E changed
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("E23").Value = "  -0.44%  "

# Row 24: E changed
$ws.Range("E24").Value = "  -1.34%  "

# Row 25: D, E changed
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.58"
$ws.Range("E25").Value = "  -0.25%  "

# Row 26: D, E changed
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.08"
$ws.Range("E26").Value = "  -0.90%  "

# Row 27: D, E changed
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.21"
$ws.Range("E27").Value = "  +0.05%  "

# Row 28: E changed
$ws.Range("E28").Value = "  +0.30%  "

# Row 29: E changed
$ws.Range("E29").Value = "  +0.56%  "

# Row 30: E changed
$ws.Range("E30").Value = "  -0.83%  "

# Row 31: E changed
$ws.Range("E31").Value = "  -0.27%  "

# Row 32: E changed
$ws.Range("E32").Value = "  +0.94%  "

# Row 33: D, E changed
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.68"
$ws.Range("E33").Value = "  +1.60%  "

# Row 34: D, E changed
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -2.76%  "

# Row 35: D, E changed
$ws.Range("D35").Value = "1.441.43"
$ws.Range("E35").Value = "  +3.05%  "

# Row 36: D, E changed
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").Value = "  +3.75%  "

# Row 37: D, E changed
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.647"
$ws.Range("E37").Value = "  -0.90%  "

# Row 38: D, E changed
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0189"
$ws.Range("E38").Value = "  +0.65%  "

# Row 39: D, E changed
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.03"
$ws.Range("E39").Value = "  -0.93%  "

# Row 40: E changed
$ws.Range("E40").Value = "  +1.93%  "

# Row 41: D, E changed
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "80.52"
$ws.Range("E41").Value = "  +0.77%  "

# Row 42: D, E changed
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  +1.01%  "

# Row 43: D, E changed
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.911"
$ws.Range("E43").Value = "  -0.85%  "

# Row 44: E changed
$ws.Range("E44").Value = "  +0.04%  "

# Row 45: E changed
$ws.Range("E45").Value = "  +2.47%  "

# Row 46: B, C, D, E changed
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.08"
$ws.Range("E46").Value = "  +0.20%  "

# Row 47: B, C, D, E changed
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.03"
$ws.Range("E47").Value = "  -1.55%  "

# Row 48: D, E changed
$ws.Range("D48").Value = "1.936.84"
$ws.Range("E48").Value = "  -0.57%  "

# Row 49: E changed
$ws.Range("E49").Value = "  +0.35%  "

# Row 50: B, C, D, E changed
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0131"
$ws.Range("E50").Value = "  -6.26%  "

# Row 51: B, C, D, E changed
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.89"
$ws.Range("E51").Value = "  -3.44%  "
